$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.354.60'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.422.57'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.50'
$ws.Range("E5").Value = '  -1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.66'
$ws.Range("E6").Value = '  +1.33%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.425.90'
$ws.Range("E8").Value = '  -1.03%  '
$ws.Range("E9").Value = '  +8.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.30'
$ws.Range("E10").Value = '  -4.08%  '
$ws.Range("E11").Value = '  +0.26%  '
$ws.Range("E12").Value = '  -0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.016.84'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000191'
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.93'
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.407.19'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.391.79'
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("E19").Value = '  -2.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.12'
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '382.50'
$ws.Range("E21").Value = '  -3.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.10'
$ws.Range("E22").Value = '  -5.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.74'
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.540'
$ws.Range("E25").Value = '  -0.92%  '
$ws.Range("E26").Value = '  +9.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.51'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("E29").Value = '  -0.12%  '
$ws.Range("E30").Value = '  +4.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.40'
$ws.Range("E31").Value = '  +0.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.03'
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.53'
$ws.Range("E33").Value = '  -2.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.44'
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.05'
$ws.Range("E36").Value = '  +2.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.14'
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.48'
$ws.Range("E38").Value = '  -1.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.026.64'
$ws.Range("E39").Value = '  +5.85%  '
$ws.Range("E40").Value = '  -0.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0760'
$ws.Range("E41").Value = '  -3.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.15'
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.50'
$ws.Range("E43").Value = '  +1.97%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0316'
$ws.Range("E44").Value = '  -1.87%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.50'
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("E46").Value = '  -1.94%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.37'
$ws.Range("E47").Value = '  +7.34%  '
$ws.Range("E48").Value = '  -2.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.871'
$ws.Range("E49").Value = '  +4.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.60'
$ws.Range("E51").Value = '  +1.86%  '
